# WR_89708709_WeekEnding_072725.xlsx edit
# Enforce single work request / single day per Excel file:
#  - refresh the "Report Generated On" timestamp
#  - update summary totals (Total Billed Amount / Total Line Items)
#  - clear the (no longer applicable) Scope ID value
#  - fill in real pricing for the remaining Thursday line items
#  - drop the duplicate "Point 06 / Inst" line item row
#  - drop the entire Friday section (this file now covers Thursday only)
#  - roll the new totals up into the TOTAL row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / summary updates -------------------------------------------------
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 09:59 AM"
$ws.Range("C8").Value = 898.72
$ws.Range("C9").Value = 2
$ws.Range("G10").Value = ""

# --- Line item pricing for the surviving Thursday rows ------------------------
$ws.Range("H16").Value = 28.6
$ws.Range("H17").Value = 870.12

# --- Remove the extra "Point 06 / Inst" duplicate line (old row 18) -----------
$ws.Rows("18").Delete()

# --- Remove the trailing blank rows + the whole Friday block (old rows 20-25,
#     now rows 19-24 after the row-18 delete above shifted everything up) -----
$ws.Rows("19:24").Delete()

# --- Roll the remaining line items into the TOTAL row (now row 18) ------------
$ws.Range("H18").Value = 898.72

Write-Host "edit complete"
